$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: "Labelname" header plus two data rows, mirroring column K's formatting
$ws.Range("L7").Value = "Labelname"
$ws.Range("L8").Value = "Exi800 - 1"
$ws.Range("L9").Value = "Exi800 - 2"

$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)

$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").Select()
